$wb = $excel.ActiveWorkbook
Write-Host ($excel.Top())
Write-Host ($excel.Left())
Write-Host ($excel.Width())
Write-Host ($excel.Height())
$excel.Top = 1170
$excel.Left = 1170
$excel.Width = 21600
$excel.Height = 13740
